# Auto-generated edit script applying the cryptos.xlsx diff
# (GitHub Actions crypto-price refresh, Fri Nov 10 14:58:40 UTC 2023)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '37.104.32'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -2.10%  '
# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.079.59'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +5.46%  '
# Row 4
$ws.Range("E4").Value = '  +0.05%  '
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '249.74'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.24%  '
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.648'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -7.14%  '
# Row 7
$ws.Range("E7").Value = '  +0.20%  '
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '50.61'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +3.93%  '
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '60.42'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.87%  '
# Row 10
$ws.Range("E10").Value = '  -4.94%  '
# Row 11
$ws.Range("E11").Value = '  -5.09%  '
# Row 12
$ws.Range("E12").Value = '  +4.22%  '
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '15.06'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -4.80%  '
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.382.60'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +6.43%  '
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.825'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.68%  '
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.090.88'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +6.28%  '
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.05'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -3.64%  '
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '36.900.10'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -2.50%  '
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '71.92'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -5.34%  '
# Row 20
$ws.Range("E20").Value = '  -5.93%  '
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.15'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -5.18%  '
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '237.87'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -7.15%  '
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.18'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.88%  '
# Row 24
$ws.Range("E24").Value = '  -0.13%  '
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.45'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -3.08%  '
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '168.78'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.02%  '
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.30'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +3.13%  '
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '20.64'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +7.55%  '
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.98'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -7.81%  '
# Row 30
$ws.Range("E30").Value = '  -6.58%  '
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.06'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +17.23%  '
# Row 32
$ws.Range("E32").Value = '  -4.13%  '
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0601'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.37%  '
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '19.93'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.12%  '
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0888'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -4.27%  '
# Row 36
$ws.Range("E36").Value = '  -0.19%  '
# Row 37
$ws.Range("B37").Value = 'LidoDAOToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.26'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +10.23%  '
# Row 38
$ws.Range("B38").Value = 'WEMIXToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.82'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.99%  '
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '4.04'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -8.01%  '
# Row 40
$ws.Range("E40").Value = '  -11.15%  '
# Row 41
$ws.Range("B41").Value = 'InjectiveProtocol'
$ws.Range("C41").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '17.49'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.61%  '
# Row 42
$ws.Range("B42").Value = 'VeChain'
$ws.Range("C42").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0221'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.85%  '
# Row 43
$ws.Range("E43").Value = '  +0.42%  '
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '97.19'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -7.95%  '
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.77'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -5.03%  '
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0876'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.94%  '
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.96'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +4.28%  '
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.300.26'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -5.45%  '
# Row 49
$ws.Range("E49").Value = '  +5.40%  '
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.256.75'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +5.97%  '
# Row 51
$ws.Range("E51").Value = '  -8.76%  '
